# Update the "Förändrad" date column (C) for all data rows (2-216)
# from 45190 (2023-09-21) to 45192 (2023-09-23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C216").Value = 45192
